$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<aral>"
$ws.Range("C2").Value = 16

# Row 3
$ws.Range("B3").Value = "<see>"
$ws.Range("C3").Value = 18

# Row 4
$ws.Range("C4").Value = 15

# Row 5
$ws.Range("C5").Value = 14

# Row 6
$ws.Range("B6").Value = "<sentence>"
$ws.Range("C6").Value = 10

# Row 8
$ws.Range("C8").Value = 16

# Row 9
$ws.Range("C9").Value = 8

# Row 10
$ws.Range("B10").Value = "<many>"
$ws.Range("C10").Value = 8

# Row 11
$ws.Range("C11").Value = 12

# Row 12
$ws.Range("C12").Value = 11

# Row 13
$ws.Range("C13").Value = 16

# Row 14
$ws.Range("C14").Value = 7

# Row 16
$ws.Range("C16").Value = 7

# Row 17
$ws.Range("C17").Value = 14

# Row 18
$ws.Range("C18").Value = 16
